# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
# (price updates + row re-ordering for the crypto listing table) described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as literal text (e.g. "249.55") in the source data.
# Force a Text number format first so Excel does not reinterpret these numeric-looking
# strings as actual numbers when the .Value is assigned below.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D24", "D25", "D28", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D50")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "249.55"
$ws.Range("D3").Value = "22.84"
$ws.Range("D4").Value = "5.406"
$ws.Range("D5").Value = "0.05623"
$ws.Range("D6").Value = "3.439"
$ws.Range("D7").Value = "6.348"
$ws.Range("D8").Value = "0.8157"
$ws.Range("D9").Value = "0.9194"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1434"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07515"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03194"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03088"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09330"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.556"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001602"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04719"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005760"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "0.006388"
$ws.Range("D20").Value = "0.005069"
$ws.Range("D21").Value = "0.001034"
$ws.Range("D24").Value = "2.162"
$ws.Range("D25").Value = "0.3296"
$ws.Range("E27").Value = "26AAXTokenAAB"
$ws.Range("D28").Value = "0.0003000"
$ws.Range("D40").Value = "0.04008"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002761"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.002917"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "0.007690"
$ws.Range("D45").Value = "0.00005572"
$ws.Range("D48").Value = "0.6752"
$ws.Range("D49").Value = "0.2019"
$ws.Range("D50").Value = "0.00002100"
